$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.977.44'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '1.636.24'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  -0.46%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.39'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5085'
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('E8').Value = '  -0.72%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06342'
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07760'
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.266'
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('D13').Value = '1.640.18'
$ws.Range('E13').Value = '  -1.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5423'
$ws.Range('E14').Value = '  -0.41%  '
$ws.Range('D15').Value = '0.0₅7692'
$ws.Range('E15').Value = '  -2.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.87'
$ws.Range('E16').Value = '  -1.70%  '
$ws.Range('D17').Value = '25.988.50'
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9998'
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '198.91'
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.407'
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('E21').Value = '  -0.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.044'
$ws.Range('E22').Value = '  +0.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  -0.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.883'
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.27'
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1202'
$ws.Range('E26').Value = '  +4.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.816'
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.57'
$ws.Range('E28').Value = '  -1.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.232'
$ws.Range('E29').Value = '  -0.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04905'
$ws.Range('E30').Value = '  -1.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.252'
$ws.Range('E31').Value = '  -0.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.168'
$ws.Range('E32').Value = '  -1.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.526'
$ws.Range('E33').Value = '  -0.55%  '
$ws.Range('E34').Value = '  -0.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9084'
$ws.Range('E35').Value = '  +1.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.581'
$ws.Range('E36').Value = '  -1.24%  '
$ws.Range('D37').Value = '1.129.64'
$ws.Range('E37').Value = '  -1.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5449'
$ws.Range('E38').Value = '  -1.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01562'
$ws.Range('E39').Value = '  -0.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9998'
$ws.Range('E40').Value = '  -0.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.522'
$ws.Range('E41').Value = '  -1.59%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8099'
$ws.Range('E42').Value = '  -1.78%  '
$ws.Range('B43').Value = 'BabyDogeCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D43').Value = '0.0₈125'
$ws.Range('E43').Value = '  +4.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '98.98'
$ws.Range('E44').Value = '  -1.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.411'
$ws.Range('E45').Value = '  -4.94%  '
$ws.Range('D46').Value = '1.776.95'
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4525'
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.87'
$ws.Range('E49').Value = '  -1.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05110'
$ws.Range('E50').Value = '  +0.99%  '
$ws.Range('E51').Value = '  -0.30%  '
